# Apply "chore: init work 31" edit to the Estudos workbook.
#
# Summary of the change:
#  - A new entry was logged for day 44865 (row 40) as an "ESTÁGIO" day, and
#    the existing "Alpha EdTeck" typo was corrected to "Alpha EdTech"
#    throughout the sheet (ASSUNTO/PRODUÇÃO columns of the Tabela1 table).
#  - Row 39 (44864) got its missing HORA F / DESCANSO / ASSUNTO filled in,
#    and its PRODUÇÃO note was changed.
#  - Row 40 (44865) got its HORA I, ASSUNTO and PRODUÇÃO filled in.
#  - The active selection moved from H39 to H38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we are working on the right sheet (there is only one: "Estudos").
$ws = $wb.Worksheets.Item("Estudos")
$ws.Activate()

# --- Row 40 (44865): add the new "ESTÁGIO + " note before the global
# rename below, so that it ends up inserted into the shared-strings table
# right after the existing strings (matching the authored order), and
# before the "Alpha EdTeck" -> "Alpha EdTech" strings get regenerated.
$ws.Range("G40").Value = "ESTÁGIO + "

# --- Global fix of the "EdTeck" typo to "EdTech" across the whole sheet.
# This touches every PRODUÇÃO cell that previously mentioned
# "Alpha EdTeck" (rows 17, 26-33, 35-38).
$ws.Cells.Replace("EdTeck", "EdTech", -4163)

# --- Row 39 (44864): fill HORA F, DESCANSO, ASSUNTO and update PRODUÇÃO.
$ws.Range("C39").Value = 0.5
$ws.Range("E39").Value = [Double]1 / 24
$ws.Range("G39").Value = "HARD"
$ws.Range("H39").Value = "Atividade voluntária no Alpha EdTech + Estudando typscript"

# --- Row 40 (44865): fill HORA I and update PRODUÇÃO (ASSUNTO already set
# above as "ESTÁGIO + ").
$ws.Range("B40").Value = [Double]8 / 24
$ws.Range("H40").Value = "Atividade voluntária no Alpha EdTech + Estágio"

# --- Move the active selection from H39 to H38.
$ws.Range("H38").Select()
